# Updates cryptocurrency price (column D) and 1h-volume-change (column E) values
# to match the latest scrape, preserving each cell's original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextCell 2 4 "70.228.82"
Set-TextCell 2 5 "  +1.24%  "
Set-TextCell 3 4 "3.443.28"
Set-TextCell 3 5 "  +1.58%  "
Set-TextCell 4 5 "  -0.02%  "
Set-TextCell 5 4 "584.45"
Set-TextCell 5 5 "  -0.53%  "
Set-TextCell 6 4 "178.69"
Set-TextCell 6 5 "  -0.64%  "
Set-TextCell 7 4 "0.600"
Set-TextCell 7 5 "  +0.75%  "
Set-TextCell 8 4 "3.433.00"
Set-TextCell 8 5 "  +1.41%  "
Set-TextCell 10 4 "0.207"
Set-TextCell 10 5 "  +6.20%  "
Set-TextCell 11 5 "  -0.29%  "
Set-TextCell 12 4 "49.05"
Set-TextCell 12 5 "  +1.19%  "
Set-TextCell 13 5 "  +1.70%  "
Set-TextCell 14 4 "693.29"
Set-TextCell 14 5 "  +2.22%  "
Set-TextCell 15 4 "3.989.59"
Set-TextCell 15 5 "  +1.39%  "
Set-TextCell 16 4 "8.70"
Set-TextCell 16 5 "  +0.99%  "
Set-TextCell 17 4 "70.112.64"
Set-TextCell 18 4 "3.436.18"
Set-TextCell 18 5 "  +1.64%  "
Set-TextCell 20 4 "17.76"
Set-TextCell 20 5 "  +0.68%  "
Set-TextCell 21 4 "11.51"
Set-TextCell 21 5 "  +2.40%  "
Set-TextCell 22 4 "0.904"
Set-TextCell 23 4 "5.52"
Set-TextCell 23 5 "  +1.80%  "
Set-TextCell 24 4 "17.13"
Set-TextCell 24 5 "  +0.12%  "
Set-TextCell 25 4 "101.41"
Set-TextCell 25 5 "  -1.92%  "
Set-TextCell 26 5 "  +0.57%  "
Set-TextCell 27 5 "  -1.39%  "
Set-TextCell 28 4 "9.65"
Set-TextCell 28 5 "  -0.50%  "
Set-TextCell 29 4 "33.68"
Set-TextCell 29 5 "  -1.33%  "
Set-TextCell 30 4 "8.80"
Set-TextCell 30 5 "  +1.15%  "
Set-TextCell 31 5 "  +3.22%  "
Set-TextCell 32 4 "3.91"
Set-TextCell 32 5 "  +8.61%  "
Set-TextCell 33 4 "574.07"
Set-TextCell 33 5 "  +3.19%  "
Set-TextCell 34 4 "11.07"
Set-TextCell 34 5 "  -0.72%  "
Set-TextCell 35 4 "58.73"
Set-TextCell 35 5 "  +1.08%  "
Set-TextCell 36 5 "  -2.57%  "
Set-TextCell 37 5 "  +0.02%  "
Set-TextCell 38 4 "3.586.19"
Set-TextCell 38 5 "  -2.67%  "
Set-TextCell 39 5 "  +0.45%  "
Set-TextCell 40 4 "35.37"
Set-TextCell 40 5 "  +1.14%  "
Set-TextCell 41 4 "0.0₃0743"
Set-TextCell 41 5 "  +6.65%  "
Set-TextCell 42 4 "3.34"
Set-TextCell 42 5 "  +1.94%  "
Set-TextCell 43 4 "2.70"
Set-TextCell 43 5 "  +0.83%  "
Set-TextCell 44 4 "3.34"
Set-TextCell 44 5 "  +1.41%  "
Set-TextCell 45 4 "0.0424"
Set-TextCell 45 5 "  +0.57%  "
Set-TextCell 46 4 "0.337"
Set-TextCell 46 5 "  -0.53%  "
Set-TextCell 47 5 "  +0.77%  "
Set-TextCell 48 4 "1.46"
Set-TextCell 48 5 "  +3.49%  "
Set-TextCell 49 4 "0.130"
Set-TextCell 49 5 "  +0.08%  "
Set-TextCell 50 4 "1.00"
Set-TextCell 50 5 "  -0.16%  "
Set-TextCell 51 4 "133.71"
Set-TextCell 51 5 "  +0.93%  "
